$d = $word.ActiveDocument

# --- 1. Merge "presents a " + "hands-on, guided workshop..." into a single run ---
$d.Content.Find.Execute(
    "presents a hands-on, guided workshop for a small group, introducing CRISPR-Cas9, a tool widely used to edit genes within organisms.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "presents a hands-on, guided workshop for a small group, introducing CRISPR-Cas9, a tool widely used to edit genes within organisms.",
    2) | Out-Null

# --- 2. Merge "...Come" + " along to learn..." into a single run ---
$d.Content.Find.Execute(
    "toolkit in laboratories. Come along to learn the theory behind this process and do an introductory bit of gene editing yourself in our ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "toolkit in laboratories. Come along to learn the theory behind this process and do an introductory bit of gene editing yourself in our ",
    2) | Out-Null

# --- 3. Merge "...as well " + "as discussing..." into a single run ---
$d.Content.Find.Execute(
    "expressed fluorescent protein, as well as discussing the applications, challenges and ethics inherent in using this technology.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "expressed fluorescent protein, as well as discussing the applications, challenges and ethics inherent in using this technology.",
    2) | Out-Null

# --- 4. Merge "For " + "adults with limited..." into a single run ---
$d.Content.Find.Execute(
    "For adults with limited formal biology training. 3 evening classes, 2 weekend days; optional one-day weekend data analysis session. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "For adults with limited formal biology training. 3 evening classes, 2 weekend days; optional one-day weekend data analysis session. ",
    2) | Out-Null

# --- 5. Merge the Mozilla Foundation paragraph runs into a single run ---
$d.Content.Find.Execute(
    "This workshop is a pilot of new educational materials, funded by the Mozilla Foundation. As such, there is no charge for participation. The trainers are not CRISPR-cas9 or sequencing experts but are keen intermediate users who wish to share the excitement and possibilities of modern biological tools. All materials are in development, so we welcome your feedback, especially if a few things are still rough around the edges!",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This workshop is a pilot of new educational materials, funded by the Mozilla Foundation. As such, there is no charge for participation. The trainers are not CRISPR-cas9 or sequencing experts but are keen intermediate users who wish to share the excitement and possibilities of modern biological tools. All materials are in development, so we welcome your feedback, especially if a few things are still rough around the edges!",
    2) | Out-Null

# --- 6. Split the "contacted by Friday 21st February." sentence into a placeholder ---
$rng = $d.Content
$rng.Find.Execute("Please complete this form. All applicants will be contacted by Friday 21st February.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    $paraStart = $rng.Start
    $paraEnd = $rng.End

    # Range for "Friday 21st February" (the part to replace with placeholder)
    $dateRng = $d.Content
    $dateRng.Find.Execute("Friday 21st February", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $dateStart = $dateRng.Start
    $dateEnd = $dateRng.End

    # Range for trailing "."
    $periodRng = $d.Range($dateEnd, $paraEnd)

    # Replace date range text with placeholder, coloured red
    $dateRng.Text = "[Organiser: Add date]"
    $dateRng.Font.Color = 255  # wdColorRed (0x0000FF = BGR 255 -> RGB FF0000)

    # Recompute the period range start after the text replacement (length differs)
    $newDateEnd = $dateStart + ("[Organiser: Add date]".Length)
    $periodRng2 = $d.Range($newDateEnd, $newDateEnd + 1)
    $periodRng2.Font.Color = 0  # wdColorBlack == RGB(0,0,0) -> themeColor text1 equivalent
    $periodRng2.Font.TextColor.ObjectThemeColor = 13  # wdThemeColorText1
}
